# Generate Report for Handback
# Replace the two source-file identifiers (UUID-based .md names) and their
# dependent hashes/timestamps across the Overview / zh-cn / de-de sheets,
# and refresh the corresponding hyperlinks' display text + targets.

$wb = $excel.ActiveWorkbook

$uuid1Old = "21a530d6-0c48-4d77-acef-69568d1e92a6"
$uuid1New = "d63d69fb-8042-45c3-aa77-f70fc98dad3a"
$uuid2Old = "662c1350-e0d6-4cc9-b0a2-2cbda34d7e9e"
$uuid2New = "ffff69cbc4f4-01dc-469f-a87a-b50bc28044f1"
$hashNew  = "db0e01fdea85202ae87950977945cd3056850bb2"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$uuid1New.md"
$wsOverview.Range("B2").Value = "e2e\$uuid1New.md"
$wsOverview.Range("G2").Value = "2016-09-02 09:13:23"

$wsOverview.Range("A3").Value = "$uuid2New.md"
$wsOverview.Range("B3").Value = "e2e\$uuid2New.md"
$wsOverview.Range("G3").Value = "2016-09-02 09:13:23"

# Rebuild the two hyperlinks on this sheet so their display text matches
# the new file names (targets follow the same github blob URL pattern).
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5ce16913ea25462b8d93e534ca7967ef9171879/e2e/$uuid1New.md", [Type]::Missing, [Type]::Missing, "e2e\$uuid1New.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5ce16913ea25462b8d93e534ca7967ef9171879/e2e/$uuid2New.md", [Type]::Missing, [Type]::Missing, "e2e\$uuid2New.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$uuid1New.$hashNew.zh-cn.xlf"

$wsZhCn.Range("A2").Value = "$uuid1New.md"
$wsZhCn.Range("G2").Value = $zhXlf
$wsZhCn.Range("H2").Value = "2016-09-02 09:13:19"
$wsZhCn.Range("I2").Value = "$uuid1New.md"
$wsZhCn.Range("J2").Value = $zhXlf
$wsZhCn.Range("K2").Value = "2016-09-02 09:13:36"

$wsZhCn.Range("A3").Value = "$uuid2New.md"
$wsZhCn.Range("G3").Value = $zhXlf
$wsZhCn.Range("H3").Value = "2016-09-02 09:13:19"
$wsZhCn.Range("I3").Value = "$uuid2New.md"
$wsZhCn.Range("J3").Value = $zhXlf
$wsZhCn.Range("K3").Value = "2016-09-02 09:13:36"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5ce16913ea25462b8d93e534ca7967ef9171879/e2e/$uuid1New.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0e59b11734f5ae1857edb48b36c21202f04cfdb4/e2e/$uuid1New.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5ce16913ea25462b8d93e534ca7967ef9171879/e2e/$uuid2New.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0e59b11734f5ae1857edb48b36c21202f04cfdb4/e2e/$uuid2New.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf = "$uuid1New.$hashNew.de-de.xlf"

$wsDeDe.Range("A2").Value = "$uuid1New.md"
$wsDeDe.Range("G2").Value = $deXlf
$wsDeDe.Range("H2").Value = "2016-09-02 09:13:23"
$wsDeDe.Range("I2").Value = "$uuid1New.md"
$wsDeDe.Range("J2").Value = $deXlf
$wsDeDe.Range("K2").Value = "2016-09-02 09:13:43"

$wsDeDe.Range("A3").Value = "$uuid2New.md"
$wsDeDe.Range("G3").Value = $deXlf
$wsDeDe.Range("H3").Value = "2016-09-02 09:13:23"
$wsDeDe.Range("I3").Value = "$uuid2New.md"
$wsDeDe.Range("J3").Value = $deXlf
$wsDeDe.Range("K3").Value = "2016-09-02 09:13:43"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5ce16913ea25462b8d93e534ca7967ef9171879/e2e/$uuid1New.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3089e3984a0358ba3f8b8499c5221e28aa6f151c/e2e/$uuid1New.md", [Type]::Missing, [Type]::Missing, "$uuid1New.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5ce16913ea25462b8d93e534ca7967ef9171879/e2e/$uuid2New.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3089e3984a0358ba3f8b8499c5221e28aa6f151c/e2e/$uuid2New.md", [Type]::Missing, [Type]::Missing, "$uuid2New.md") | Out-Null
